$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-08-16 Friday", $false, $true, $false, $false, $false, `
    $true, 1, $false, "2024-08-17 Saturday", 2) | Out-Null

# Update each arithmetic expression in the table (old values are unique in the
# document, so a global literal Find/Replace is safe for every pair).

$d.Content.Find.Execute("1+45=46", $false, $true, $false, $false, $false, $true, 1, $false, "83-72=11", 2) | Out-Null

$d.Content.Find.Execute("72+11=83", $false, $true, $false, $false, $false, $true, 1, $false, "62+13=75", 2) | Out-Null

$d.Content.Find.Execute("6+9=15", $false, $true, $false, $false, $false, $true, 1, $false, "68-33=35", 2) | Out-Null

$d.Content.Find.Execute("97-31=66", $false, $true, $false, $false, $false, $true, 1, $false, "31+12=43", 2) | Out-Null

$d.Content.Find.Execute("57-22=35", $false, $true, $false, $false, $false, $true, 1, $false, "72-51=21", 2) | Out-Null

$d.Content.Find.Execute("76-60=16", $false, $true, $false, $false, $false, $true, 1, $false, "49-4=45", 2) | Out-Null

$d.Content.Find.Execute("55-7=48", $false, $true, $false, $false, $false, $true, 1, $false, "30+3=33", 2) | Out-Null

$d.Content.Find.Execute("90-54=36", $false, $true, $false, $false, $false, $true, 1, $false, "0+29=29", 2) | Out-Null

$d.Content.Find.Execute("63+25=88", $false, $true, $false, $false, $false, $true, 1, $false, "33+28=61", 2) | Out-Null

$d.Content.Find.Execute("10+77=87", $false, $true, $false, $false, $false, $true, 1, $false, "9+67=76", 2) | Out-Null

$d.Content.Find.Execute("14+32=46", $false, $true, $false, $false, $false, $true, 1, $false, "81-24=57", 2) | Out-Null

$d.Content.Find.Execute("56-50=6", $false, $true, $false, $false, $false, $true, 1, $false, "77-71=6", 2) | Out-Null

$d.Content.Find.Execute("15+34=49", $false, $true, $false, $false, $false, $true, 1, $false, "75-25=50", 2) | Out-Null

$d.Content.Find.Execute("52+0=52", $false, $true, $false, $false, $false, $true, 1, $false, "35-4=31", 2) | Out-Null

$d.Content.Find.Execute("75-62=13", $false, $true, $false, $false, $false, $true, 1, $false, "48-43=5", 2) | Out-Null

$d.Content.Find.Execute("10+86=96", $false, $true, $false, $false, $false, $true, 1, $false, "96-5=91", 2) | Out-Null

$d.Content.Find.Execute("41+58=99", $false, $true, $false, $false, $false, $true, 1, $false, "63-18=45", 2) | Out-Null

$d.Content.Find.Execute("95-10=85", $false, $true, $false, $false, $false, $true, 1, $false, "70+17=87", 2) | Out-Null

$d.Content.Find.Execute("35+36=71", $false, $true, $false, $false, $false, $true, 1, $false, "23+36=59", 2) | Out-Null

$d.Content.Find.Execute("46-35=11", $false, $true, $false, $false, $false, $true, 1, $false, "42+55=97", 2) | Out-Null

$d.Content.Find.Execute("41-6=35", $false, $true, $false, $false, $false, $true, 1, $false, "90-70=20", 2) | Out-Null

$d.Content.Find.Execute("90+2=92", $false, $true, $false, $false, $false, $true, 1, $false, "31-8=23", 2) | Out-Null

$d.Content.Find.Execute("16-15=1", $false, $true, $false, $false, $false, $true, 1, $false, "8+75=83", 2) | Out-Null

$d.Content.Find.Execute("62-46=16", $false, $true, $false, $false, $false, $true, 1, $false, "3+75=78", 2) | Out-Null

$d.Content.Find.Execute("57+4=61", $false, $true, $false, $false, $false, $true, 1, $false, "31+59=90", 2) | Out-Null

$d.Content.Find.Execute("98-46=52", $false, $true, $false, $false, $false, $true, 1, $false, "66-56=10", 2) | Out-Null

$d.Content.Find.Execute("37-9=28", $false, $true, $false, $false, $false, $true, 1, $false, "26+34=60", 2) | Out-Null

$d.Content.Find.Execute("16-6=10", $false, $true, $false, $false, $false, $true, 1, $false, "59-44=15", 2) | Out-Null

$d.Content.Find.Execute("93-22=71", $false, $true, $false, $false, $false, $true, 1, $false, "14-1=13", 2) | Out-Null

$d.Content.Find.Execute("53+29=82", $false, $true, $false, $false, $false, $true, 1, $false, "5+25=30", 2) | Out-Null

$d.Content.Find.Execute("96-89=7", $false, $true, $false, $false, $false, $true, 1, $false, "76+10=86", 2) | Out-Null

$d.Content.Find.Execute("74-5=69", $false, $true, $false, $false, $false, $true, 1, $false, "42+4=46", 2) | Out-Null

$d.Content.Find.Execute("67+14=81", $false, $true, $false, $false, $false, $true, 1, $false, "76+19=95", 2) | Out-Null

$d.Content.Find.Execute("99-17=82", $false, $true, $false, $false, $false, $true, 1, $false, "28+36=64", 2) | Out-Null

$d.Content.Find.Execute("52-46=6", $false, $true, $false, $false, $false, $true, 1, $false, "53+13=66", 2) | Out-Null

$d.Content.Find.Execute("93-30=63", $false, $true, $false, $false, $false, $true, 1, $false, "75-62=13", 2) | Out-Null

$d.Content.Find.Execute("60+9=69", $false, $true, $false, $false, $false, $true, 1, $false, "7+30=37", 2) | Out-Null

$d.Content.Find.Execute("28+14=42", $false, $true, $false, $false, $false, $true, 1, $false, "19-3=16", 2) | Out-Null

$d.Content.Find.Execute("99-54=45", $false, $true, $false, $false, $false, $true, 1, $false, "74-43=31", 2) | Out-Null

$d.Content.Find.Execute("96-67=29", $false, $true, $false, $false, $false, $true, 1, $false, "62-32=30", 2) | Out-Null

$d.Content.Find.Execute("15-14=1", $false, $true, $false, $false, $false, $true, 1, $false, "3+96=99", 2) | Out-Null

$d.Content.Find.Execute("36+4=40", $false, $true, $false, $false, $false, $true, 1, $false, "43+48=91", 2) | Out-Null

$d.Content.Find.Execute("87-74=13", $false, $true, $false, $false, $false, $true, 1, $false, "79-1=78", 2) | Out-Null

$d.Content.Find.Execute("35+5=40", $false, $true, $false, $false, $false, $true, 1, $false, "37+59=96", 2) | Out-Null

$d.Content.Find.Execute("65-11=54", $false, $true, $false, $false, $false, $true, 1, $false, "41-21=20", 2) | Out-Null

$d.Content.Find.Execute("27-21=6", $false, $true, $false, $false, $false, $true, 1, $false, "97-69=28", 2) | Out-Null

$d.Content.Find.Execute("24+12=36", $false, $true, $false, $false, $false, $true, 1, $false, "39+51=90", 2) | Out-Null

$d.Content.Find.Execute("82+12=94", $false, $true, $false, $false, $false, $true, 1, $false, "35+49=84", 2) | Out-Null

$d.Content.Find.Execute("72+1=73", $false, $true, $false, $false, $false, $true, 1, $false, "89-75=14", 2) | Out-Null

$d.Content.Find.Execute("13+39=52", $false, $true, $false, $false, $false, $true, 1, $false, "84-75=9", 2) | Out-Null

$d.Content.Find.Execute("45+45=90", $false, $true, $false, $false, $false, $true, 1, $false, "73-66=7", 2) | Out-Null

$d.Content.Find.Execute("16+19=35", $false, $true, $false, $false, $false, $true, 1, $false, "2+73=75", 2) | Out-Null

$d.Content.Find.Execute("37-30=7", $false, $true, $false, $false, $false, $true, 1, $false, "90-70=20", 2) | Out-Null

$d.Content.Find.Execute("53-15=38", $false, $true, $false, $false, $false, $true, 1, $false, "93-80=13", 2) | Out-Null

$d.Content.Find.Execute("78-61=17", $false, $true, $false, $false, $false, $true, 1, $false, "14+36=50", 2) | Out-Null

$d.Content.Find.Execute("49+1=50", $false, $true, $false, $false, $false, $true, 1, $false, "6+35=41", 2) | Out-Null

$d.Content.Find.Execute("16+58=74", $false, $true, $false, $false, $false, $true, 1, $false, "9+68=77", 2) | Out-Null

$d.Content.Find.Execute("70-49=21", $false, $true, $false, $false, $false, $true, 1, $false, "78+20=98", 2) | Out-Null

$d.Content.Find.Execute("11+72=83", $false, $true, $false, $false, $false, $true, 1, $false, "89-46=43", 2) | Out-Null

$d.Content.Find.Execute("36+23=59", $false, $true, $false, $false, $false, $true, 1, $false, "19+32=51", 2) | Out-Null

$d.Content.Find.Execute("67-55=12", $false, $true, $false, $false, $false, $true, 1, $false, "53+39=92", 2) | Out-Null

$d.Content.Find.Execute("20+13=33", $false, $true, $false, $false, $false, $true, 1, $false, "61+36=97", 2) | Out-Null

$d.Content.Find.Execute("85+7=92", $false, $true, $false, $false, $false, $true, 1, $false, "49-13=36", 2) | Out-Null

$d.Content.Find.Execute("69-48=21", $false, $true, $false, $false, $false, $true, 1, $false, "80-48=32", 2) | Out-Null

$d.Content.Find.Execute("92-80=12", $false, $true, $false, $false, $false, $true, 1, $false, "61-12=49", 2) | Out-Null

$d.Content.Find.Execute("82-4=78", $false, $true, $false, $false, $false, $true, 1, $false, "14+64=78", 2) | Out-Null

$d.Content.Find.Execute("22+55=77", $false, $true, $false, $false, $false, $true, 1, $false, "15+63=78", 2) | Out-Null

$d.Content.Find.Execute("94-35=59", $false, $true, $false, $false, $false, $true, 1, $false, "22+5=27", 2) | Out-Null

$d.Content.Find.Execute("73-53=20", $false, $true, $false, $false, $false, $true, 1, $false, "81-36=45", 2) | Out-Null

$d.Content.Find.Execute("52+22=74", $false, $true, $false, $false, $false, $true, 1, $false, "21-13=8", 2) | Out-Null

$d.Content.Find.Execute("13+20=33", $false, $true, $false, $false, $false, $true, 1, $false, "10+32=42", 2) | Out-Null

$d.Content.Find.Execute("32+8=40", $false, $true, $false, $false, $false, $true, 1, $false, "89+4=93", 2) | Out-Null

$d.Content.Find.Execute("74-52=22", $false, $true, $false, $false, $false, $true, 1, $false, "17+35=52", 2) | Out-Null

$d.Content.Find.Execute("17-3=14", $false, $true, $false, $false, $false, $true, 1, $false, "99-36=63", 2) | Out-Null

$d.Content.Find.Execute("2+11=13", $false, $true, $false, $false, $false, $true, 1, $false, "64-50=14", 2) | Out-Null

$d.Content.Find.Execute("77+13=90", $false, $true, $false, $false, $false, $true, 1, $false, "71+28=99", 2) | Out-Null

$d.Content.Find.Execute("39+60=99", $false, $true, $false, $false, $false, $true, 1, $false, "97-83=14", 2) | Out-Null

$d.Content.Find.Execute("69+24=93", $false, $true, $false, $false, $false, $true, 1, $false, "58-44=14", 2) | Out-Null

$d.Content.Find.Execute("88-1=87", $false, $true, $false, $false, $false, $true, 1, $false, "83-31=52", 2) | Out-Null

$d.Content.Find.Execute("75+9=84", $false, $true, $false, $false, $false, $true, 1, $false, "9+61=70", 2) | Out-Null

$d.Content.Find.Execute("26+4=30", $false, $true, $false, $false, $false, $true, 1, $false, "75-66=9", 2) | Out-Null

$d.Content.Find.Execute("20+18=38", $false, $true, $false, $false, $false, $true, 1, $false, "61-53=8", 2) | Out-Null

$d.Content.Find.Execute("74+10=84", $false, $true, $false, $false, $false, $true, 1, $false, "67-11=56", 2) | Out-Null

$d.Content.Find.Execute("88-48=40", $false, $true, $false, $false, $false, $true, 1, $false, "30+7=37", 2) | Out-Null

$d.Content.Find.Execute("75-50=25", $false, $true, $false, $false, $false, $true, 1, $false, "88-6=82", 2) | Out-Null

$d.Content.Find.Execute("48-34=14", $false, $true, $false, $false, $false, $true, 1, $false, "20+16=36", 2) | Out-Null

$d.Content.Find.Execute("77-30=47", $false, $true, $false, $false, $false, $true, 1, $false, "70+3=73", 2) | Out-Null

$d.Content.Find.Execute("85+2=87", $false, $true, $false, $false, $false, $true, 1, $false, "56-33=23", 2) | Out-Null

$d.Content.Find.Execute("0+53=53", $false, $true, $false, $false, $false, $true, 1, $false, "71-67=4", 2) | Out-Null

$d.Content.Find.Execute("85-40=45", $false, $true, $false, $false, $false, $true, 1, $false, "5+12=17", 2) | Out-Null

$d.Content.Find.Execute("97-0=97", $false, $true, $false, $false, $false, $true, 1, $false, "17+23=40", 2) | Out-Null

$d.Content.Find.Execute("58+20=78", $false, $true, $false, $false, $false, $true, 1, $false, "12+50=62", 2) | Out-Null

$d.Content.Find.Execute("7+5=12", $false, $true, $false, $false, $false, $true, 1, $false, "25+53=78", 2) | Out-Null

$d.Content.Find.Execute("58-18=40", $false, $true, $false, $false, $false, $true, 1, $false, "22-19=3", 2) | Out-Null

$d.Content.Find.Execute("94-3=91", $false, $true, $false, $false, $false, $true, 1, $false, "45+23=68", 2) | Out-Null

$d.Content.Find.Execute("47+40=87", $false, $true, $false, $false, $false, $true, 1, $false, "61-34=27", 2) | Out-Null

$d.Content.Find.Execute("25+37=62", $false, $true, $false, $false, $false, $true, 1, $false, "57-17=40", 2) | Out-Null

$d.Content.Find.Execute("73-64=9", $false, $true, $false, $false, $false, $true, 1, $false, "33+20=53", 2) | Out-Null

$d.Content.Find.Execute("31+4=35", $false, $true, $false, $false, $false, $true, 1, $false, "46-41=5", 2) | Out-Null

$d.Content.Find.Execute("61-0=61", $false, $true, $false, $false, $false, $true, 1, $false, "13+18=31", 2) | Out-Null
